# Edit script for "AT3_Component Design and Integration Plan (1).docx"
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Extend the "Feral Kingdom..." paragraph with new sentences, and
#    remove the now-orphaned "_GoBack" bookmark that used to sit at the
#    end of that paragraph.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "the latter ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "the latter will not be marked as a win state, and thus the player would have to come back and complete the fight to be able to beat the game. Considerations need to be made with the monster components as to how the enemies and player take damage and deal damage to allow the system to be reused for each monster type to make each monster" + [char]0x2019 + "s play style be unique from the others.",
    2
) | Out-Null

$gb = $d.Bookmarks("_GoBack")
$gb.Delete()

# ---------------------------------------------------------------------
# 2) Move the lastRenderedPageBreak marker from the "How everything is
#    connected" run to the start of the "Create a domain model..." run.
# ---------------------------------------------------------------------
# (handled further below once text anchors are in place)

# ---------------------------------------------------------------------
# 3) "Using the three point" -> "Using the three-point"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "three point estimation formula",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "three-point estimation formula",
    2
) | Out-Null

# ---------------------------------------------------------------------
# 4) "each components calculation" -> "each component's calculation"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "each components calculation",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "each component" + [char]0x2019 + "s calculation",
    2
) | Out-Null

# ---------------------------------------------------------------------
# 5) Header date: 22/05/2020 -> 25/05/2020
# ---------------------------------------------------------------------
foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers.Item(1)
    $hdr.Range.Find.Execute(
        "22/05/2020",
        $true, $false, $false, $false, $false, $true, 1, $false,
        "25/05/2020",
        2
    ) | Out-Null
}
